$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.289.47"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.871.37"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7082"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.66"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07774"
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3091"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.03"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08395"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.864.85"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.241"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7112"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.03"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "29.297.08"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.067"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008175"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.98"
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "2.118.71"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.751"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1588"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.22"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.016"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.45"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.504"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.399"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.289"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.302"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05327"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.935"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7447"
$ws.Range("E36").Value = "  -6.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "1.230.90"
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.724"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.557"
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.83"
$ws.Range("E42").Value = "  +6.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8848"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "2.015.71"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5192"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.793"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.399"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4306"
$ws.Range("E51").Value = "  +0.28%  "
